# ---------------------------------------------------------------------------
# Applies the "DMEM layout" sheet addition + register-convention rework
# described by the commit diff to the "reg convention.xlsx" workbook.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("DMEM state ")

# ===========================================================================
# 1. Sheet1 ("Register Convention--Logo") rework
# ===========================================================================

# -- Row 4: new note in H4 --------------------------------------------------
$ws1.Range("H4").Value = "Color code: 0, 1, 2, 3, etc…"

# -- Row 5: new mini-header for the Orientation/Color legend (M5:O5) -------
$ws1.Range("M5").Value = "Orientation Code"
$ws1.Range("M5:N5").HorizontalAlignment = 7   # xlCenterContinuous
$ws1.Range("O5").Value = "Color"
$ws1.Range("O5").HorizontalAlignment = -4108  # xlCenter

# -- Rows 6-9: move the old North/East/South/West + Red/Green/Blue/White ---
# -- legend (cols G:H) over to the new cols N:O, and replace G:H with the --
# -- new Const/Direction labels ---------------------------------------------
$ws1.Range("G6:H9").Copy($ws1.Range("N6:O9"))

$ws1.Range("M6").Value = 0
$ws1.Range("M7").Value = 1
$ws1.Range("M8").Value = 2
$ws1.Range("M9").Value = 3

$ws1.Range("G6").Value = "Const"
$ws1.Range("H6").ClearContents()

$ws1.Range("F7").Value = "reserved temp reg"
$ws1.Range("G7").Value = "Direction"
$ws1.Range("H7").ClearContents()

$ws1.Range("F8").ClearContents()
$ws1.Range("G8").ClearContents()
$ws1.Range("H8").ClearContents()

$ws1.Range("F9").ClearContents()
$ws1.Range("G9").ClearContents()
$ws1.Range("H9").ClearContents()

$ws1.Activate()
$ws1.Range("I12").Select()

# ===========================================================================
# 2. "DMEM state " sheet - just a selection change
# ===========================================================================
$ws2.Activate()
$ws2.Range("H5").Select()

# ===========================================================================
# 3. New "DMEM layout" sheet
# ===========================================================================
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "DMEM layout"

$ws3.Range("F3").Value = "Addr"

$ws3.Range("F4").Value = 0
$ws3.Range("H4").Value = "Function pointer"

$ws3.Range("C5").Value = "Turtle Images * 4"

$ws3.Range("F6").Value = 900
$ws3.Range("F8").Value = 1800
$ws3.Range("F10").Value = 2700

$ws3.Range("C13").Value = "States"
$ws3.Range("F13").Value = 3600

$ws3.Range("C21").Value = "*Memory grows downward"
$ws3.Range("F21").Value = 4000

# Thick box borders around C4:E12 and C13:E22 (two "cards")
$ws3.Range("C4:E12").BorderAround(1, 4)   # xlContinuous, xlThick
$ws3.Range("C13:E22").BorderAround(1, 4)  # xlContinuous, xlThick

$ws3.Range("G23").Select()

$ws1.Activate()
$ws1.Range("I12").Select()
